$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow edits, re-protect afterwards
$ws.Unprotect()

# Update the confidential disclaimer date (2021-03-23 -> 2021-03-24)
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-03-24 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15
$ws.Range("D2").Value = 0.05519512553617553
$ws.Range("E2").Value = -0.008283501617488476

$ws.Range("D3").Value = 0.02309446079499865
$ws.Range("E3").Value = -0.004305043050430557

$ws.Range("D4").Value = 0.0315129213325535
$ws.Range("E4").Value = -0.02337260961947074

$ws.Range("D5").Value = 0.03109832212716095
$ws.Range("E5").Value = 0.02024793388429758

$ws.Range("D6").Value = 0.03225351860809159
$ws.Range("E6").Value = 0.003522898842476252

$ws.Range("D7").Value = 0.01840493121183709
$ws.Range("E7").Value = 0.004506775597559187

$ws.Range("D8").Value = 0.004758083866128257
$ws.Range("E8").Value = -0.0447761194029852

$ws.Range("D9").Value = 0.0064753226630152
$ws.Range("E9").Value = -0.001671192813870848

$ws.Range("D10").Value = 0.06963643423689628
$ws.Range("E10").Value = 0.002913752913753065

$ws.Range("D11").Value = 0.06979875692742518
$ws.Range("E11").Value = 0.002906976744186052

$ws.Range("D12").Value = 0.1487849781387859
$ws.Range("E12").Value = 0.005382209615244493

$ws.Range("D13").Value = 0.3930887435193005
$ws.Range("E13").Value = 0.001316251316251238

$ws.Range("D14").Value = 0.1158984010376316
$ws.Range("E14").Value = 0.003921568627450966

$ws.Range("E15").Value = 0.001487714512479199

# Re-apply protection to match original workbook state
$ws.Protect()
